$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (new Price text, new Volume(1h) text); $null means "leave unchanged"
$updates = @(
    @{ Row = 2; D = "26.403.03"; E = "  -3.24%  " }
    @{ Row = 3; D = "1.801.53"; E = "  -2.84%  " }
    @{ Row = 4; D = "1.008"; E = "  +0.58%  " }
    @{ Row = 5; D = $null; E = "  +0.55%  " }
    @{ Row = 6; D = "308.13"; E = "  -2.04%  " }
    @{ Row = 7; D = "0.4523"; E = "  -1.70%  " }
    @{ Row = 8; D = "0.3641"; E = "  -1.85%  " }
    @{ Row = 9; D = "0.07067"; E = "  -3.13%  " }
    @{ Row = 10; D = "0.8678"; E = "  -2.43%  " }
    @{ Row = 11; D = "0.07762"; E = "  -0.91%  " }
    @{ Row = 12; D = "19.23"; E = "  -4.42%  " }
    @{ Row = 13; D = "1.821.33"; E = "  -0.83%  " }
    @{ Row = 14; D = "5.244"; E = "  -2.76%  " }
    @{ Row = 15; D = "6.312"; E = "  -3.27%  " }
    @{ Row = 16; D = "85.95"; E = "  -6.02%  " }
    @{ Row = 17; D = "1.009"; E = "  +0.56%  " }
    @{ Row = 18; D = "0.000008535"; E = "  -4.46%  " }
    @{ Row = 19; D = "1.008"; E = "  +0.56%  " }
    @{ Row = 20; D = "26.458.88"; E = "  -3.12%  " }
    @{ Row = 21; D = "14.17"; E = "  -3.94%  " }
    @{ Row = 22; D = "4.950"; E = "  -3.21%  " }
    @{ Row = 23; D = $null; E = "  -1.51%  " }
    @{ Row = 24; D = "1.974"; E = "  +2.88%  " }
    @{ Row = 25; D = "150.14"; E = "  -1.16%  " }
    @{ Row = 26; D = "17.85"; E = "  -3.32%  " }
    @{ Row = 27; D = "1.980"; E = "  -3.93%  " }
    @{ Row = 28; D = "112.64"; E = "  -2.90%  " }
    @{ Row = 29; D = "4.851"; E = "  -4.33%  " }
    @{ Row = 30; D = "0.08631"; E = "  -2.19%  " }
    @{ Row = 31; D = "3.025"; E = "  -2.23%  " }
    @{ Row = 32; D = "0.7261"; E = "  -6.06%  " }
    @{ Row = 33; D = "4.429"; E = "  -1.93%  " }
    @{ Row = 34; D = "1.109"; E = "  -5.26%  " }
    @{ Row = 35; D = "1.007"; E = "  +0.56%  " }
    @{ Row = 36; D = "2.530"; E = "  -8.40%  " }
    @{ Row = 37; D = "1.075"; E = "  -0.36%  " }
    @{ Row = 38; D = "0.01916"; E = "  -1.94%  " }
    @{ Row = 39; D = $null; E = "  -2.69%  " }
    @{ Row = 40; D = "0.05049"; E = "  -4.12%  " }
    @{ Row = 41; D = "6.943"; E = "  -1.49%  " }
    @{ Row = 42; D = "0.4890"; E = "  -4.58%  " }
    @{ Row = 43; D = "0.1565"; E = "  -4.45%  " }
    @{ Row = 44; D = "8.093"; E = "  -3.70%  " }
    @{ Row = 45; D = "1.008"; E = "  +0.60%  " }
    @{ Row = 46; D = "0.4587"; E = "  -4.29%  " }
    @{ Row = 47; D = "9.937"; E = "  -4.62%  " }
    @{ Row = 48; D = "100.90"; E = "  -1.51%  " }
    @{ Row = 49; D = "1.577"; E = "  -4.05%  " }
    @{ Row = 50; D = "0.05978"; E = "  -3.86%  " }
    @{ Row = 51; D = "63.34"; E = "  -3.70%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Range("D" + $u.Row)
        # Force text storage so numeric-looking price strings (e.g. "1.008")
        # are not reinterpreted as numbers - they must stay strings, matching
        # the source inlineStr cells.
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
